$wb = $excel.ActiveWorkbook

# 1) Duplicate "Planilha15" into a new sheet placed at the end, renamed "sheet nova"
$srcSheet = $wb.Worksheets.Item("Planilha15")
$srcSheet.Copy([System.Type]::Missing, $srcSheet)
$newSheet = $wb.Worksheets.Item($srcSheet.Index + 1)
$newSheet.Name = "sheet nova"

# 2) Fix Planilha15: the title-row filler cells (B1:E1, H1:K1, F1, L1) should no longer
#    carry the bold/centered "code" style - reset them to plain default formatting.
$ws = $wb.Worksheets.Item("Planilha15")

$rng1 = $ws.Range("B1:E1")
$rng1.Font.Name = "Calibri"
$rng1.Font.Bold = $false
$rng1.HorizontalAlignment = 1
$rng1.VerticalAlignment = -4107

$rng2 = $ws.Range("F1")
$rng2.Font.Name = "Calibri"
$rng2.Font.Bold = $false
$rng2.HorizontalAlignment = 1
$rng2.VerticalAlignment = -4107

$rng3 = $ws.Range("H1:K1")
$rng3.Font.Name = "Calibri"
$rng3.Font.Bold = $false
$rng3.HorizontalAlignment = 1
$rng3.VerticalAlignment = -4107

$rng4 = $ws.Range("L1")
$rng4.Font.Name = "Calibri"
$rng4.Font.Bold = $false
$rng4.HorizontalAlignment = 1
$rng4.VerticalAlignment = -4107

$wb.Worksheets.Item($wb.Worksheets.Count).Select()
